$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = -0.5934017483515905
$ws.Range("J5").Value = 0.4529499003166974
$ws.Range("K5").Value = 0.2367433508890373
$ws.Range("L5").Value = 2.628059292783453
